$wb = $excel.ActiveWorkbook

# --- remove Sheet3 (album data moves into renamed Sheet2 instead) ---
$wb.Worksheets.Item("Sheet3").Delete()

# --- rename remaining sheets ---
$wb.Worksheets.Item("Sheet1").Name = "T_TRACK_TYPE_MST"
$wb.Worksheets.Item("Sheet2").Name = "T_ALBUM_TYPE_MST"

$ws = $wb.Worksheets.Item("T_ALBUM_TYPE_MST")

# --- populate new strings in the exact order they were first typed so the
#     shared-string table grows in the same sequence as the authored file ---
$ws.Range("C3").Value = "片尾曲"

$ws.Range("C2").Value = "片头曲"
$r1 = $ws.Range("C2").Characters(2,1)
$r1.Font.Name = "ＭＳ Ｐゴシック"
$r2 = $ws.Range("C2").Characters(3,1)
$r2.Font.Name = "ＭＳ Ｐゴシック"

$ws.Range("B7").Value = "OST"

$ws.Range("C7").Value = "原声音乐"
$r3 = $ws.Range("C7").Characters(4,1)
$r3.Font.Name = "ＭＳ Ｐゴシック"

$ws.Range("C9").Value = "广播剧(角色)"
$r4 = $ws.Range("C9").Characters(3,5)
$r4.Font.Name = "ＭＳ Ｐゴシック"

$ws.Range("C8").Value = "短声音集"

$ws.Range("A1").Value = "ALBUM_TYPE_ID"
$ws.Range("B1").Value = "ALBUM_TYPE_NAME"

# --- remaining cells reuse already-known shared strings / are plain numbers ---
$ws.Range("C1").Value = "DESCRIPTION"

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "OP"

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "ED"

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "IN"
$ws.Range("C4").Value = "插入曲"

$ws.Range("A5").Value = 13
$ws.Range("B5").Value = "IM"
$ws.Range("C5").Value = "印象曲"

$ws.Range("A6").Value = 14
$ws.Range("B6").Value = "CS"
$ws.Range("C6").Value = "角色曲"

$ws.Range("A7").Value = 20

$ws.Range("A8").Value = 30
$ws.Range("B8").Value = "Sound"

$ws.Range("A9").Value = 40
$ws.Range("B9").Value = "Drama"

$ws.Range("A10").Value = 50
$ws.Range("B10").Value = "Radio"
$ws.Range("C10").Value = "广播节目"
